$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 18 - this shifts the existing rows 18-55 down
# to rows 19-56 (and the sheet dimension grows from A1:R55 to A1:R56).
$ws.Rows.Item(18).Insert()

# Populate the newly-inserted row 18 with a new weekly price observation for
# "Agrícola del Norte S.A. de Arica" / Albahaca, matching the surrounding
# rows' static columns and carrying new date/price figures.
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44980
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100112052
$ws.Range("G18").Value = "Albahaca"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 900
$ws.Range("N18").Value = "$/paquete"
$ws.Range("O18").Value = "Región de Arica y Parinacota"
$ws.Range("P18").Value = 900
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
